# "Danh sach can bo nhan vien" - update API new code
# Replace the staff table data (rows 4-13) with a new set of 9 staff records
# (rows 4-12), so the Excel Table (Table1) shrinks from A3:I13 to A3:I12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old 10th data row entirely. Because this is inside the Table
# (ListObject), deleting the worksheet row also resizes Table1 / its
# AutoFilter from A3:I13 down to A3:I12 and fixes up the sheet dimension.
$ws.Rows("13").Delete()

# Row 4
$ws.Range("B4").Value = "NV72573840"
$ws.Range("C4").Value = "Đinh Yên Nam"
$ws.Range("D4").Value = "'074232234"
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = "Hóa"
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = "x"

# Row 5
$ws.Range("B5").Value = "NV92831013"
$ws.Range("C5").Value = "Vũ Kiến Nam"
$ws.Range("D5").Value = "'044292284"
$ws.Range("E5").Value = "Tổ Sử Địa "
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = "x"

# Row 6
$ws.Range("B6").Value = "NV99829582"
$ws.Range("C6").Value = "Lê Ngọc Nam"
$ws.Range("D6").Value = "'087342212"
$ws.Range("E6").Value = "Tổ Hóa - Sinh"
$ws.Range("F6").Value = "Toán"
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = "x"

# Row 7
$ws.Range("B7").Value = "NV59507985"
$ws.Range("C7").Value = "Nguyễn Mạnh Nam"
$ws.Range("D7").Value = "'069309572"
$ws.Range("E7").Value = "Tổ Hóa - Sinh"
$ws.Range("F7").Value = ""
$ws.Range("G7").Value = ""
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = "x"

# Row 8
$ws.Range("B8").Value = "NV66662587"
$ws.Range("C8").Value = "Lê Như Nam"
$ws.Range("D8").Value = "'031561639"
$ws.Range("E8").Value = "Tổ Hóa - Sinh"
$ws.Range("F8").Value = "Văn"
$ws.Range("G8").Value = ""
$ws.Range("H8").Value = ""
$ws.Range("I8").Value = "x"

# Row 9
$ws.Range("B9").Value = "NV18993783"
$ws.Range("C9").Value = "Nguyễn Minh Nam"
$ws.Range("D9").Value = "'054327222"
$ws.Range("E9").Value = "Tổ Lý"
$ws.Range("F9").Value = ""
$ws.Range("G9").Value = ""
$ws.Range("H9").Value = ""
$ws.Range("I9").Value = "x"

# Row 10
$ws.Range("B10").Value = "NV89320025"
$ws.Range("C10").Value = "Hồ Phương Nam"
$ws.Range("D10").Value = "'023439844"
$ws.Range("E10").Value = "Tổ Sử Địa "
$ws.Range("F10").Value = ""
$ws.Range("G10").Value = ""
$ws.Range("H10").Value = ""
$ws.Range("I10").Value = "x"

# Row 11
$ws.Range("B11").Value = "NV52112691"
$ws.Range("C11").Value = "Lê Đona Nam"
$ws.Range("D11").Value = "'034396549"
$ws.Range("E11").Value = "Tổ Lý"
$ws.Range("F11").Value = "Toán"
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = ""
$ws.Range("I11").Value = "x"

# Row 12
$ws.Range("B12").Value = "NV89033488"
$ws.Range("C12").Value = "Vũ Gia Nam"
$ws.Range("D12").Value = "'087350901"
$ws.Range("E12").Value = "Tổ Lý"
$ws.Range("F12").Value = "Văn"
$ws.Range("G12").Value = ""
$ws.Range("H12").Value = ""
$ws.Range("I12").Value = "x"
